$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "cervical-cancer-screening" row (old row 11); everything
# below shifts up by one.
$ws.Rows(11).Delete()

# Old row 10's B value pointed at the (now removed) cervical-cancer URL;
# repoint it at the ncicdr url so row 9/10 form a loop pair.
$ws.Range("B10").Value = "http://www.dev03.webmd.com/cancer/tc/ncicdr0000593694-families"

# New "third hop" column with relationship tags for a few rows. Shared
# strings are appended in first-use order, so touch them in the same order
# they appear in the target workbook (duplicate, loop, multihop).
$ws.Range("C8").Value = "duplicate"
$ws.Range("C11").Value = "duplicate"
$ws.Range("C9").Value = "loop"
$ws.Range("C10").Value = "loop"
$ws.Range("C5").Value = "multihop"
$ws.Range("C6").Value = "multihop"

# New hyperlinks on A9 and B10 (Excel auto-creates the "Hyperlink" cell
# style/font the first time a hyperlink is added to the workbook).
$ws.Hyperlinks.Add($ws.Range("A9"), "http://www.dev03.webmd.com/cancer/tc/ncicdr0000593694-families")
$ws.Hyperlinks.Add($ws.Range("B10"), "http://www.dev03.webmd.com/cancer/tc/pediatric-supportive-care-pdq-supportive-care---patient-information-nci-families")

# Narrower first column.
$ws.Columns(1).ColumnWidth = 117.14

# Update the view: select C10 and make sure the window is scrolled back to A.
$ws.Range("C10").Select()
